# Update "江西-漫展信息" workbook per commit "Update gh-pages to output generated at 456a3b4"
#
# Changes:
#   1. Sheet "展览" (Exhibitions): refresh the "想去人数" (interest count) column
#      for ten existing events.
#   2. Sheet "演出" (Performances): append a new event row
#      (南昌·《梁祝》65周年大型交响音乐会-风兔子交响乐团, 2024-09-16).
#   3. Sheet "全部类型" (All types, a combined/sorted view of every event):
#      insert the same new event in date order and refresh the same
#      "想去人数" counts for the events that also live on this sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. 展览 — simple cell updates
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    "F2"  = 152
    "F3"  = 1802
    "F4"  = 39
    "F9"  = 559
    "F12" = 84
    "F18" = 5125
    "F19" = 59
    "F20" = 838
    "F22" = 2284
    "F25" = 2131
}
foreach ($addr in $expoUpdates.Keys) {
    $wsExpo.Range($addr).Value = $expoUpdates[$addr]
}

# ---------------------------------------------------------------------
# 2. 演出 — insert the new row at the bottom (row 4)
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Rows.Item(4).Insert()

# Copy the number-format/border/bold style used by column A on the row
# above, so the new index cell matches the rest of the table.
$wsShow.Range("A3").Copy() | Out-Null
$wsShow.Range("A4").PasteSpecial(-4122) | Out-Null

$wsShow.Range("A4").Value = 3

# "开始时间" is stored as plain text ("2024-09-16"), not a real date — force
# text formatting first so Excel doesn't silently convert it to a date
# serial, then drop back to the default style so no stray number format
# is left attached to the cell.
$wsShow.Range("B4").NumberFormat = "@"
$wsShow.Range("B4").Value = "2024-09-16"
$wsShow.Range("B4").Style = "Normal"

$wsShow.Range("C4").Value = "南昌·《梁祝》65周年大型交响音乐会-风兔子交响乐团"
$wsShow.Range("D4").Value = "象湖新城东祥路昌南文化中心五号馆 昌南文化中心大剧院"
$wsShow.Range("E4").Value = "2024.09.16 19:30-09.16 20:50"
$wsShow.Range("F4").Value = 0
$wsShow.Range("G4").Value = 140
$wsShow.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=90515"
$wsShow.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202408/muKn0Ygv1723107475651.jpeg"

# ---------------------------------------------------------------------
# 3. 全部类型 — insert the same new event in date order (row 26, pushing
#    the two following rows down to 27/28), then refresh the counts.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Rows.Item(26).Insert()

$wsAll.Range("A25").Copy() | Out-Null
$wsAll.Range("A26").PasteSpecial(-4122) | Out-Null

$wsAll.Range("A26").Value = 25

# Column A is a plain literal running index (not a formula), so the rows
# pushed down by the insert keep their old numbers and need to be bumped
# by one by hand: old row26 (index 25) is now row27, old row27 (index 26)
# is now row28.
$wsAll.Range("A27").Value = 26
$wsAll.Range("A28").Value = 27

$wsAll.Range("B26").NumberFormat = "@"
$wsAll.Range("B26").Value = "2024-09-16"
$wsAll.Range("B26").Style = "Normal"

$wsAll.Range("C26").Value = "南昌·《梁祝》65周年大型交响音乐会-风兔子交响乐团"
$wsAll.Range("D26").Value = "象湖新城东祥路昌南文化中心五号馆 昌南文化中心大剧院"
$wsAll.Range("E26").Value = "2024.09.16 19:30-09.16 20:50"
$wsAll.Range("F26").Value = 0
$wsAll.Range("G26").Value = 140
$wsAll.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=90515"
$wsAll.Range("I26").Value = "//i2.hdslb.com/bfs/openplatform/202408/muKn0Ygv1723107475651.jpeg"

$allUpdates = @{
    "F2"  = 152
    "F3"  = 1802
    "F4"  = 39
    "F9"  = 559
    "F12" = 84
    "F18" = 5125
    "F20" = 59
    "F22" = 838
    "F24" = 2284
    "F28" = 2131
}
foreach ($addr in $allUpdates.Keys) {
    $wsAll.Range($addr).Value = $allUpdates[$addr]
}
